$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder speaker_variant id/label pairs (column B/C) across rows 2-20,
# and clear the is_prefered flag (column D) for every data row, matching
# an export with no is_pref and no lev distance.

$ws.Range("B2").Value = '#israel'
$ws.Range("C2").Value = 'Israel'
$ws.Range("D2").Value = ""

$ws.Range("B3").Value = '#israel.'
$ws.Range("C3").Value = 'Israel.'
$ws.Range("D3").Value = ""

$ws.Range("B4").Value = '#aristobulus'
$ws.Range("C4").Value = 'Aristobulus'
$ws.Range("D4").Value = ""

$ws.Range("B5").Value = '#amassai'
$ws.Range("C5").Value = 'Amassai'
$ws.Range("D5").Value = ""

$ws.Range("B6").Value = '#cogitatio'
$ws.Range("C6").Value = 'Cogitatio'
$ws.Range("D6").Value = ""

$ws.Range("B7").Value = '#iohanna'
$ws.Range("C7").Value = 'Iohanna'
$ws.Range("D7").Value = ""

$ws.Range("B8").Value = '#cognitio'
$ws.Range("C8").Value = 'Cognitio'
$ws.Range("D8").Value = ""

$ws.Range("B9").Value = '#amal-marma-te-zamen'
$ws.Range("C9").Value = 'Amal Marma te zamen'
$ws.Range("D9").Value = ""

$ws.Range("B10").Value = '#amal.       marma'
$ws.Range("C10").Value = 'Amal.       Marma'
$ws.Range("D10").Value = ""

$ws.Range("B11").Value = '#eubulus'
$ws.Range("C11").Value = 'Eubulus'
$ws.Range("D11").Value = ""

$ws.Range("B12").Value = '#achazib'
$ws.Range("C12").Value = 'Achazib'
$ws.Range("D12").Value = ""

$ws.Range("B13").Value = '#bode'
$ws.Range("C13").Value = 'Bode'
$ws.Range("D13").Value = ""

$ws.Range("B14").Value = '#amal.       marma.'
$ws.Range("C14").Value = 'Amal.       Marma.'
$ws.Range("D14").Value = ""

$ws.Range("B15").Value = '#aristobolus'
$ws.Range("C15").Value = 'Aristobolus'
$ws.Range("D15").Value = ""

$ws.Range("B16").Value = '#precatio'
$ws.Range("C16").Value = 'Precatio'
$ws.Range("D16").Value = ""

$ws.Range("B17").Value = '#amal'
$ws.Range("C17").Value = 'Amal'
$ws.Range("D17").Value = ""

$ws.Range("B18").Value = '#neregel'
$ws.Range("C18").Value = 'Neregel'
$ws.Range("D18").Value = ""

$ws.Range("B19").Value = '#marma'
$ws.Range("C19").Value = 'Marma'
$ws.Range("D19").Value = ""

$ws.Range("B20").Value = '#demophon'
$ws.Range("C20").Value = 'Demophon'
$ws.Range("D20").Value = ""

Write-Host "applied reorder + cleared is_prefered column"
